$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The long "encoded" strings used in column B (kept verbatim from the
# original workbook, referenced here so rows 6/7 can be rebuilt).
$long1 = "U2FsdGVkX18CeviQ24hIEERMC5lllAkCyTq8qB+K1J8ELLtppxyKsUFxOE+nntq9/WW38YXJQgmprt6+xWNAmDhQSQC6lExckmwooXEjTNnJ3+9TvopXVeVD4S0/efSymJkoILmDy9RhiqXoSqODsFf1mw71a6OLivpxueg/q+qqYdvq9yTI5iNQnz8y5ZvghwX8Tnm229QpdNMizeUnveiSMqY7iXdMtRUQvaeCOZiQmAcsaONQgiZHeeIJfpeKgFkOVDkheWibIj5j8MfxlLI/fFNKAEvF/SJnCDktvpPdSOUz0PHB2+E7GdBSsjK5"

# Replace the last two remaining time-slot rows (previously rows 6 & 7,
# which held "12:25-12:30" / "12:50-12:55") with the new values.
$ws.Range("C6").Value = "18:55-19:0"
$ws.Range("B7").Value = $long1
$ws.Range("C7").Value = "19:0-19:5"

# Remove the trailing rows 8-11 (old entries for 14:10-14:15, 14:15-14:20,
# 22:45-22:50, 22:50-22:55) so only rows 1-7 remain.
$ws.Range("A8:C11").EntireRow.Delete()

# Update the selection to match the saved view state.
$ws.Range("B13").Select()
